# Actualización automática hashcode
# Updates the hashcode values (column B) for a set of keyed rows in the
# "hashcode" metadata sheet, matching the author's regeneration of hashes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = "fef132f1f5ff59d4d500645fdae2eafc"
$ws.Range("B51").Value = "6375b0c9e19540fcced85b59935248e4"
$ws.Range("B54").Value = "efcf55780469d9ac7c006c2c7d748dee"
$ws.Range("B80").Value = "92d169bd7e3ee99ff67be31a0999254a"
$ws.Range("B108").Value = "e50b6189d81b10ac5e23f5b6a4c25788"
$ws.Range("B159").Value = "86a32b40bf3869218dbb5318ac73dde7"
$ws.Range("B169").Value = "75ad2a5365ea8a72ca5ddbbc28b828fb"
$ws.Range("B227").Value = "4115b0982b41c5732ea5b747d4dcb9e7"
$ws.Range("B232").Value = "9219d792f0111c25326012c6094d7a13"
$ws.Range("B339").Value = "885d675495acea9740f1c7bb31cfbbaa"
$ws.Range("B420").Value = "930e9bd628ccd09c643cd2b4a4b8cfad"
$ws.Range("B464").Value = "c64fea71094245a6b65dbbf602a9480b"
$ws.Range("B483").Value = "894927cd864488d2c75750887fe2af0f"
$ws.Range("B506").Value = "8c961637837e75f1424ec97ae6e05c47"
$ws.Range("B508").Value = "f4ecf7d3761c99fd246bf4d08bdd9a00"
$ws.Range("B524").Value = "4eb260a2ce0bb392501a2b7815ef8433"
$ws.Range("B548").Value = "2c1cd70e120f1618c514ea20d26acb33"
$ws.Range("B555").Value = "781565fc03d4b8852605f066d47696e9"
$ws.Range("B582").Value = "6207f2a46a8039f5c6d33709bcefc05b"
$ws.Range("B583").Value = "321013199fdf99fd35c8b704b3092c4e"
$ws.Range("B624").Value = "98d74cdd8f1992c38d3de5c4f237d050"
$ws.Range("B635").Value = "d91be6043d4519e7a2106349ed286d2a"
$ws.Range("B673").Value = "003bd1a3349afac2db993828b457c703"
$ws.Range("B688").Value = "3ec944572790bd9c3345656754a008a7"
$ws.Range("B693").Value = "9fb4f12c6301f274e92ec36147b7ab74"
$ws.Range("B711").Value = "1ff4be5db5422c6642b16c8b24afbeb3"
$ws.Range("B741").Value = "4cd8d12abb0ad061a5045bafd15a0c72"
$ws.Range("B827").Value = "7d618c8bf09746d171da3abda4a9112e"
$ws.Range("B858").Value = "1c3639405e3a421b3fe54b1923dd3333"
$ws.Range("B882").Value = "c9c849f03081bb7a17b5eba5feebb7ea"
